# "Make it work with Handwriting"
# Append 11 new rows (91-101) of handwriting-recognised data to Sheet1, and
# trim the stray fractional-day time component off D90.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- D90: drop the time-of-day fraction, keep just the date serial -------
$ws.Cells.Item(90, 4).Value = 44936

# --- Helper: write a piece of (possibly multi-line / CRLF) OCR'd quantity
# text into a column-C cell without Excel's auto-number-detection mangling
# it (column C's default format is a plain decimal number format). We
# briefly force the cell to Text ("@") format while assigning the literal
# string, then switch the cell back to the sheet's usual decimal format and
# turn wrapping on so the multi-line capture shows fully (matches the new
# "numFmtId 2 + wrapText" style used for these OCR'd cells).
function Set-OcrQty($cell, [string]$text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.NumberFormat = "0.00"
    $cell.WrapText = $true
}

# Row 91 - plain continuation row (clean values, like existing rows)
$ws.Cells.Item(91, 1).Value = 45270
$ws.Cells.Item(91, 2).Value = " SALATA "
$ws.Cells.Item(91, 3).Value = 500
$ws.Cells.Item(91, 4).Value = 44937

# Row 92 - quantity captured with a stray CRLF from handwriting OCR
$ws.Cells.Item(92, 1).Value = 45270
$ws.Cells.Item(92, 2).Value = " SALATA "
Set-OcrQty $ws.Cells.Item(92, 3) " 500`r`n"
$ws.Cells.Item(92, 4).Value = 44937
$ws.Rows.Item(92).RowHeight = 30

# Row 93 (column A/B keep the sheet's usual formats - the OCR'd date/produs
# text just doesn't parse cleanly as a date/number, so Excel stores it as
# text without needing any format juggling)
$ws.Cells.Item(93, 1).Value = "13-11 "
$ws.Cells.Item(93, 2).Value = " LOBODA "
Set-OcrQty $ws.Cells.Item(93, 3) " 300`r`n"
$ws.Cells.Item(93, 4).Value = 44937
$ws.Rows.Item(93).RowHeight = 30

# Row 94
$ws.Cells.Item(94, 1).Value = "14-11 "
$ws.Cells.Item(94, 2).Value = " CEAPA "
Set-OcrQty $ws.Cells.Item(94, 3) " 100`r`n15-11 "
$ws.Cells.Item(94, 4).Value = 44937
$ws.Rows.Item(94).RowHeight = 30

# Row 95
$ws.Cells.Item(95, 1).Value = " 16-11 "
$ws.Cells.Item(95, 2).Value = " USTUROI "
Set-OcrQty $ws.Cells.Item(95, 3) " 10`r`n"
$ws.Cells.Item(95, 4).Value = 44937
$ws.Rows.Item(95).RowHeight = 30

# Row 96 - back to clean values
$ws.Cells.Item(96, 1).Value = "17-11 "
$ws.Cells.Item(96, 2).Value = " DOVLEAC "
$ws.Cells.Item(96, 3).Value = 5
$ws.Cells.Item(96, 4).Value = 44937

# Rows 97-101 repeat the same sequence with a later "Data Adaugare" stamp
$ws.Cells.Item(97, 1).Value = 45270
$ws.Cells.Item(97, 2).Value = " SALATA "
Set-OcrQty $ws.Cells.Item(97, 3) " 500`r`n"
$ws.Cells.Item(97, 4).Value = 44937.84
$ws.Rows.Item(97).RowHeight = 30

$ws.Cells.Item(98, 1).Value = "13-11 "
$ws.Cells.Item(98, 2).Value = " LOBODA "
Set-OcrQty $ws.Cells.Item(98, 3) " 300`r`n"
$ws.Cells.Item(98, 4).Value = 44937.84
$ws.Rows.Item(98).RowHeight = 30

$ws.Cells.Item(99, 1).Value = "14-11 "
$ws.Cells.Item(99, 2).Value = " CEAPA "
Set-OcrQty $ws.Cells.Item(99, 3) " 100`r`n15-11 "
$ws.Cells.Item(99, 4).Value = 44937.84
$ws.Rows.Item(99).RowHeight = 30

$ws.Cells.Item(100, 1).Value = " 16-11 "
$ws.Cells.Item(100, 2).Value = " USTUROI "
Set-OcrQty $ws.Cells.Item(100, 3) " 10`r`n"
$ws.Cells.Item(100, 4).Value = 44937.84
$ws.Rows.Item(100).RowHeight = 30

$ws.Cells.Item(101, 1).Value = "17-11 "
$ws.Cells.Item(101, 2).Value = " DOVLEAC "
$ws.Cells.Item(101, 3).Value = 5
$ws.Cells.Item(101, 4).Value = 44937.84

# --- Restore the view so the newly appended rows are visible -------------
$ws.Range("D89").Select()
